$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v1 = @'
51.You use the HTML5 API to determine the current location of the user.
'@
$ws.Range("A54").Value = $v1

$v2 = @'
if (navigator.geolocation) {
    navigator.geolocation.getCurrentPosition(function (position) {
        var latitude = position.coords.latitude;
        var longitude = position.coords.longitude;
        // query database to find the nearest petrol station.
    });
}
'@
$ws.Range("B54").Value = $v2

$v3 = @'
My answer is correct.Answer explanation:The HTML Geolocation API introduced the geolocation property of the navigator object.To check if the browser is compaitable with the Geolocation API you can use if(navigator.geolocation) {}.The window.location object refers to the current page address (URL) instead of the geographical location of the user/device.To get the current position you should call the navigator.geolocation.getCurrentPosition() method and a position object is returned in the callback.The coords attribute contains the geographical coordinates which include the latitude, longitude, altitude, speed, etc.
'@
$ws.Range("D54").Value = $v3

$v4 = @'
52.You are using the Geolocation application programming interface (API) to determine a user's location.You need to retrieve the latitude and longitude coordinates.
'@
$ws.Range("A55").Value = $v4

$v5 = @'
window.navigator.geolocation.getCurrentPosition(function (position) {
    var latitude = position.coords.latitude;
    var longitude = position.coords.longitude;
});
'@
$ws.Range("B55").Value = $v5

$v6 = @'
My answer is correct.Answer explanation: The Geolocation API starts with the geolocation object.You can reference this object by accessing the geolocation property of the navigator object which represents the client as a browser.The getCurrentPosition function of the geolocation object accepts a callback function as its first parameter.This function is called if the geographic position is successfully obtained.The signature of the callback function must contain one parameter that represents a position object.The coords property of the position object returns a reference to a coordinates object that specifies the geographic information about the position.The coordinates object contains two properties named latitude and longitude that represent the position's latitude and langitude respectively.
'@
$ws.Range("D55").Value = $v6

$v7 = @'
53.A web page needs to store data for its users.The storage mechanism must meet these requirements: *It must allow data storage of up to 1MB. *The data must be retrievable after the user closes and reopens the browser. *The data must not be transmited to the Web server.
'@
$ws.Range("A56").Value = $v7

$v8 = @'
Local Storage
'@
$ws.Range("B56").Value = $v8

$v9 = @'
My answer is correct. Answer explanation: You should use the local storage.Local storage allows you to store up to 10MB of data on the client.The data is not automatically transmitted to the Web server with each browser request.When you close and reopen the browser the data is still retrievable.You access local storage through the localStorage property of the window object.
'@
$ws.Range("D56").Value = $v9

$v10 = @'
53.You are creating web site that uses the Application Cache API.You must ensure that resources named Main.js and Main.css are cached.You must ensure that Functions.js is never cached.
'@
$ws.Range("A57").Value = $v10

$v11 = @'
CACHE MANIFEST
CACHE:
Main.js
Main.css
FALLBACK:
Functions.js
'@
$ws.Range("C57").Value = $v11

$v12 = @'
CACHE MANIFEST
Main.js
Main.css
NETWORK:
Functions.js
'@
$ws.Range("B57").Value = $v12

$v13 = @'
My answer is incorrect. Answer explanation: The first line in a cache manifest file should specify the words CACHE MANIFEST.Cache manifest files are then divided into one or more of the following sections: CACHE, FALLBACK and NETWORK.Each section name must be placed on a single line and end with a colon.The CACHE section specifies the resources that should be cached.The FALLBACK section covers section specifies the resources that should be used if resources cannot be downloaded to be cached.The NETWORK section specifies resoures that should never be cached.If you do not specify a section name by default the CACHE section is assumed.In this scenario the Main.js and Main.css resources appear without a section header.Therefore they are associated with the CACHE section indicating that they should be cached.The Functions.js resource is specified in the NETWORK section indicating that it should never be cached.
'@
$ws.Range("D57").Value = $v13

$v14 = @'
54.А web page needs to save the value of an <input> element.The storage mechanism must meet these requirements: *It must allow data storage of up to 1MB *The data must be retrievable after the user closes and reopens the browser. *The data must not be transmitted to the Web server.
'@
$ws.Range("A58").Value = $v14

$v15 = @'
localStorage.userdata = $("#userData").val();
'@
$ws.Range("B58").Value = $v15

$v16 = @'
My answer is correct.Answer explanation: This code uses local storage.Local storage allows you to store up to 10MB of data on the client.The data is not automatically transmitted to the web server with each browser request.When you close and reopen the browser the data is still retrievable.You access local storage through the localStorage property of the window object.
'@
$ws.Range("D58").Value = $v16

$v17 = @'
55.You create a web site that uses AppCache.You want to always cache resources named App.js and App.css locally.You must ensure that Test.js is never cached.
'@
$ws.Range("A59").Value = $v17

$v18 = @'
CACHE MANIFEST:
CACHE:
App.js
App.css
NETWORK:
Test.js
'@
$ws.Range("B59").Value = $v18

$v19 = @'
My answer is correct.Answer explanation:In this scenario you should place the App.js and App.css resources within the CACHE section indicating that they should be cached.You should place the Test.js resource in the NETWORK section indicating that it should NEVER BE CACHED.
'@
$ws.Range("D59").Value = $v19

$v20 = @'
57.You received a JSON string from a web service.To consume the data you to convert the JSON string to JS object.
'@
$ws.Range("A60").Value = $v20

$v21 = @'
1. JSON.parse 2. jQuery.parseJSON
'@
$ws.Range("B60").Value = $v21

$v22 = @'
JSON.stringify
'@
$ws.Range("C60").Value = $v22

$v23 = @'
My answers are incorrect.Answers explanation:1. JSON.parse is a correct option because it deserializes JSON text to a JS object. 2. jQuery.parseJSON is also a correct option because it converts JSON text to a JS object.
'@
$ws.Range("D60").Value = $v23

$v24 = @'
58.You want to asynchronously load a plain text file using XMLHttpRequest.Which are the correct steps?
'@
$ws.Range("A61").Value = $v24

$v25 = @'
1. Assign a new instance of XMLHttpRequest to xhr 2. Register a handler for event onreadystatechange 3. Invoke xhr.open("GET", url, true) 4. Invoke xhr.send()
'@
$ws.Range("B61").Value = $v25

$v26 = @'
My answers are correct. Answers explanation: 1.To use the XMLHttpRequest object you must first create a new instance of it. 2.You should register an event handler to its onreadystatechange event in order to process the response data. 3.The open method should be called next.It accepts three parameters: HTTP method, URL and isAsync.You should use "GET" method to load a plain text file.The isAsync parameter should be set to true for asynchronous requests. 4.The send method should be called last to initiate the request.
'@
$ws.Range("D61").Value = $v26

$v27 = @'
59.You write the following code to retrive the title of exam 70-480:                                                                       var request = new XMLHttpRequest();
request.open("GET", "http://service.measureup.com/Exams", false);
request.send();
var attribute = doc.selectSingleNode("//Exam[@ID='70-480']/@Title");
var title = attribute.value;                                                                                                   You need to add a code line at the empty space.
'@
$ws.Range("A62").Value = $v27

$v28 = @'
var doc = request.responseXML;
'@
$ws.Range("B62").Value = $v28

$v29 = @'
My answer is correct. Answer explanation: The responseXML property returns an IXMLDOMDocument object that represents the loaded XML data.This object contains a selectSingleNode method that allows you to search for an element, attribute or text node.
'@
$ws.Range("D62").Value = $v29

$v30 = @'
60.You write the following code to retrieve the title of the first exam:                                                                                     var request = new XMLHttpRequest();
request.open("GET", "http://service.measureup.com/Exams", false);
request.send();
var title = doc.Exams[0].Title;                                                                                                                               You need to add code line at the empty space.
'@
$ws.Range("A63").Value = $v30

$v31 = @'
var doc = JSON.parse(request.responseText);
'@
$ws.Range("B63").Value = $v31

$v32 = @'
My answer is correct.Answer explanation: The responseText property returns a string that represents the data retrieved.In this scenario the data is a JSON-encoded string.This code uses the parse function of the JSON object to convert the string data to an object.
'@
$ws.Range("D63").Value = $v32

$v33 = @'
61.You need to use the IXMLDOMDocument object to determine the total cost of a sandwich and lemonade.
'@
$ws.Range("A64").Value = $v33

$v34 = @'
var totalPrice = 0;
var prices = doc.selectNodes("//Item[@Name='Sandwich' or @Name='Lemonade']/@Price");
for (var index = 0; index < prices.length; index++){
    totalPrice += prices[index].value;
}
'@
$ws.Range("B64").Value = $v34

$v35 = @'
var totalPrice = 0;
var prices = doc.selectNodes("//Item[Name='Sandwich' or Name='Lemonade']/Price");
for (var index = 0; index < prices.length; index++){
    totalPrice += prices[index].value;
}
'@
$ws.Range("C64").Value = $v35

$v36 = @'
My answer is incorrect. Answer explanation: This code calls the selectNode method passing to it an XPATH expression that returns all Price attributes where the Name attribute is equal to Sandwich or Lemonade.It then iterates throgh the returned attributes and adds their values.
'@
$ws.Range("D64").Value = $v36

$v37 = @'
63.The following JS code retrieves data from a web service:                                                                                      var request = new XMLHttpRequest();
request.open("GET", "GetData", true);
request.send();
var data = request.responseText;                                                                                                              After this code runs the data variable is empty.However if you add the following code alert("Testing") the data variable contains the correct value.
'@
$ws.Range("A65").Value = $v37

$v38 = @'
1. Add the following code between lines 01 and 02:                                                                                                 var data = null;
request.onreadystatechange = function () {
    if (request.readyState == 4 && request.status == 200) {
        data = request.responseText;
    }  
};                                                                                                                                                 2. Remove line 04
'@
$ws.Range("B65").Value = $v38

$v39 = @'
My answers are correct.Answers exolanation: 1. The problem in this scenario is that the web service is accessed asynchronously as indicated by the third parameter of the open function of the XMLHttpRequest object.When this parameter is set to true you must handle the onreadystatechange event which is raised after the data is available. 2. You should also remove line 04 because you must access the data in the event handler of the onreadystatechange event.
'@
$ws.Range("D65").Value = $v39

$v40 = @'
63.You retrieve the data as string named menu.You need to convert the string to an object and remove the Calories and Type properties so that they do not exist as part of the object.
'@
$ws.Range("A66").Value = $v40

$v41 = @'
var newMenu = JSON.parse(menu, function (key, value) {
    var newValue =  value;
    switch (key) {
        case "Calories":
        case "Type":
            newValue = undefined;
            break;
    }
    return newValue;
});
'@
$ws.Range("B66").Value = $v41

$v42 = @'
var newMenu = JSON.stringify(menu, function (key, value) {
    var newValue =  value;
    switch (key) {
        case "Calories":
        case "Type":
            newValue = undefined;
            break;
    }
    return newValue;
});
'@
$ws.Range("C66").Value = $v42

$v43 = @'
My answer is incorrect.Answer explanation: In this scenario the data retrieved is a JSON-encoded string.The parse function of the JSON object converts a JSON-encoded string to an object.The second parameter of the parse function is an optional callback that allows you to change the values of properties or remove properties altogheter.The callback accepts the property name as the first parameter and the property value as the second parameter.The return value of the callback represents the new value of the associated property.If you return undefined for a property the property is removed from the resulting object.
'@
$ws.Range("D66").Value = $v43

$v44 = @'
64.You retrieve the data as a string named data.You need to to convert the string to an object and remove the Responsibility and UNLOC properties so that they do not exist as part of the object.
'@
$ws.Range("A67").Value = $v44

$v45 = @'
1. JSON.parse 2. undefined
'@
$ws.Range("B67").Value = $v45

$v46 = @'
2. null
'@
$ws.Range("C67").Value = $v46

$v47 = @'
My answer is incorrect.Answer explanation:In this scenario the data retrieved is a JSON-encoded string.The parse function of the JSON onject converts a JSON-encoded string to an object.The second parameter of the parse function is an optional callback that allows you to change the values of properties or remove properties altogether.The callback accepts the property name as the first parameter and the property value as the second parameter.The return value of the callback represents the new value of the associated property.If you then return undefined for a property the property is removed from the resulting object.
'@
$ws.Range("D67").Value = $v47

$v48 = @'
65.You want to use an HTML5 input field to capture user's email address.You need to ensure that an email address is always supplied as part of the subscription request.
'@
$ws.Range("A68").Value = $v48

$v49 = @'
<input type="email" name="email" required />
'@
$ws.Range("B68").Value = $v49

$v50 = @'
My answers are correct.Answers explanation: The most appropriate HTML5 input type for email address is "email".Standarts complaint browsers will validate if the entered value is in valid email format.The "required" attribute ensures that no submission is made until a value is given.
'@
$ws.Range("D68").Value = $v50

$ws.Range("D67").Select()
